$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right column (B) 5 -> 4, Wrong column (C) -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): Right column (B) 75 -> 60, Wrong column (C) -5 -> -10
$ws.Range("B12").Value = 60
$ws.Range("C12").Value = -10

# Row 12, column E: summary text "75 / 140" -> "50 / 112"
$ws.Range("E12").Value = "50 / 112"
